$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Meter Serial No." column (A) used to store generator serial numbers either
# as raw numbers (row 2) or as shared-string hex codes without the "009-" site
# prefix (rows 3-10), rendered via a custom text number format (numFmtId 49).
# All devices should now be processed uniformly as plain text values carrying
# the "009-" prefix, using the sheet's default (General) formatting.

$serials = @{
  2  = "009-80845"
  3  = "009-80E1F"
  4  = "009-80DCD"
  5  = "009-80E2A"
  6  = "009-80E29"
  7  = "009-80B76"
  8  = "009-80B1E"
  9  = "009-80B1C"
  10 = "009-80B13"
}

# Drop the special "text" number format (numFmtId 49) from column A / its cells
# so every cell falls back to the workbook's single default style.
$ws.Columns.Item(1).ClearFormats()

foreach ($row in $serials.Keys | Sort-Object) {
  $cell = $ws.Cells.Item($row, 1)
  $cell.Value = $serials[$row]
}

# Column A is a bit wider now that it holds the "009-" prefixed codes.
$ws.Columns.Item(1).ColumnWidth = 19.92

# Selection moved from the old "next empty row" (A11) to the newly edited block.
$ws.Range("A3:A6").Select()
